# This script applies the "Updated symbol list" data refresh for cryptos.xlsx
# (prices / 1h volume changes, plus a few coin re-rankings on rows 8-17).
#
# Text-look-alike numeric strings (e.g. "303.79", "3.20%") must be written back
# as literal text (matching the original inlineStr cells), not auto-converted by
# Excel into numbers/percentages. We force that by temporarily switching the cell
# to a text NumberFormat before assigning .Value, then restoring the default style
# so no stray formatting is left behind on the cell.
function Set-TextValue($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; Col=4; Val='303.79'; IsText=$false},
    @{Row=2; Col=5; Val='3.20%'; IsText=$false},
    @{Row=3; Col=4; Val='33.92'; IsText=$false},
    @{Row=3; Col=5; Val='9.11%'; IsText=$false},
    @{Row=4; Col=4; Val='5.161'; IsText=$false},
    @{Row=4; Col=5; Val='4.75%'; IsText=$false},
    @{Row=5; Col=4; Val='0.07830'; IsText=$false},
    @{Row=5; Col=5; Val='6.16%'; IsText=$false},
    @{Row=6; Col=4; Val='2.413'; IsText=$false},
    @{Row=6; Col=5; Val='4.84%'; IsText=$false},
    @{Row=7; Col=4; Val='8.043'; IsText=$false},
    @{Row=7; Col=5; Val='4.51%'; IsText=$false},
    @{Row=8; Col=2; Val='GateToken'; IsText=$true},
    @{Row=8; Col=3; Val='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; IsText=$true},
    @{Row=8; Col=4; Val='3.907'; IsText=$false},
    @{Row=8; Col=5; Val='4.06%'; IsText=$false},
    @{Row=9; Col=2; Val='MXToken'; IsText=$true},
    @{Row=9; Col=3; Val='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; IsText=$true},
    @{Row=9; Col=4; Val='0.9348'; IsText=$false},
    @{Row=9; Col=5; Val='2.30%'; IsText=$false},
    @{Row=10; Col=2; Val='LiechtensteinCryptoassetsExchange'; IsText=$true},
    @{Row=10; Col=3; Val='https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'; IsText=$true},
    @{Row=10; Col=4; Val='0.09788'; IsText=$false},
    @{Row=10; Col=5; Val='17.01%'; IsText=$false},
    @{Row=11; Col=2; Val='WazirX'; IsText=$true},
    @{Row=11; Col=3; Val='https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; IsText=$true},
    @{Row=11; Col=4; Val='0.1778'; IsText=$false},
    @{Row=11; Col=5; Val='5.20%'; IsText=$false},
    @{Row=12; Col=2; Val='MandalaExchangeToken'; IsText=$true},
    @{Row=12; Col=3; Val='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; IsText=$true},
    @{Row=12; Col=4; Val='0.08500'; IsText=$false},
    @{Row=12; Col=5; Val='4.03%'; IsText=$false},
    @{Row=13; Col=2; Val='BitrueCoin'; IsText=$true},
    @{Row=13; Col=3; Val='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; IsText=$true},
    @{Row=13; Col=4; Val='0.03354'; IsText=$false},
    @{Row=13; Col=5; Val='7.51%'; IsText=$false},
    @{Row=14; Col=2; Val='BitMartToken'; IsText=$true},
    @{Row=14; Col=3; Val='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; IsText=$true},
    @{Row=14; Col=4; Val='0.09923'; IsText=$false},
    @{Row=14; Col=5; Val='-1.42%'; IsText=$false},
    @{Row=15; Col=2; Val='BitForexToken'; IsText=$true},
    @{Row=15; Col=3; Val='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; IsText=$true},
    @{Row=15; Col=4; Val='0.001480'; IsText=$false},
    @{Row=15; Col=5; Val='-1.97%'; IsText=$false},
    @{Row=16; Col=2; Val='TigerCash'; IsText=$true},
    @{Row=16; Col=3; Val='https://coinranking.com/coin/6hIn06L2+tigercash-tch'; IsText=$true},
    @{Row=16; Col=4; Val='0.005688'; IsText=$false},
    @{Row=16; Col=5; Val='-0.26%'; IsText=$false},
    @{Row=17; Col=2; Val='LEO'; IsText=$true},
    @{Row=17; Col=3; Val='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; IsText=$true},
    @{Row=17; Col=4; Val='3.466'; IsText=$false},
    @{Row=17; Col=5; Val='-0.49%'; IsText=$false},
    @{Row=18; Col=4; Val='2.162'; IsText=$false},
    @{Row=18; Col=5; Val='4.00%'; IsText=$false},
    @{Row=19; Col=4; Val='0.3367'; IsText=$false},
    @{Row=19; Col=5; Val='1.13%'; IsText=$false},
    @{Row=20; Col=5; Val='2.95%'; IsText=$false},
    @{Row=21; Col=4; Val='4.286'; IsText=$false},
    @{Row=21; Col=5; Val='8.03%'; IsText=$false},
    @{Row=22; Col=4; Val='0.2290'; IsText=$false},
    @{Row=22; Col=5; Val='9.09%'; IsText=$false},
    @{Row=23; Col=4; Val='0.04648'; IsText=$false},
    @{Row=23; Col=5; Val='2.30%'; IsText=$false},
    @{Row=24; Col=5; Val='0.97%'; IsText=$false},
    @{Row=25; Col=4; Val='0.004417'; IsText=$false},
    @{Row=25; Col=5; Val='1.86%'; IsText=$false},
    @{Row=26; Col=4; Val='0.0001296'; IsText=$false},
    @{Row=26; Col=5; Val='-0.34%'; IsText=$false},
    @{Row=27; Col=4; Val='0.0003385'; IsText=$false},
    @{Row=27; Col=5; Val='-0.25%'; IsText=$false},
    @{Row=39; Col=4; Val='0.01745'; IsText=$false},
    @{Row=39; Col=5; Val='8.70%'; IsText=$false},
    @{Row=40; Col=4; Val='0.04810'; IsText=$false},
    @{Row=40; Col=5; Val='8.10%'; IsText=$false},
    @{Row=41; Col=4; Val='0.007787'; IsText=$false},
    @{Row=41; Col=5; Val='6.34%'; IsText=$false},
    @{Row=42; Col=4; Val='0.009799'; IsText=$false},
    @{Row=42; Col=5; Val='10.97%'; IsText=$false},
    @{Row=43; Col=4; Val='0.1415'; IsText=$false},
    @{Row=43; Col=5; Val='6.84%'; IsText=$false},
    @{Row=44; Col=4; Val='0.002079'; IsText=$false},
    @{Row=44; Col=5; Val='0.88%'; IsText=$false},
    @{Row=45; Col=4; Val='0.009134'; IsText=$false},
    @{Row=45; Col=5; Val='-0.57%'; IsText=$false},
    @{Row=46; Col=4; Val='0.00006091'; IsText=$false},
    @{Row=46; Col=5; Val='1.52%'; IsText=$false},
    @{Row=47; Col=4; Val='0.00000000748'; IsText=$false},
    @{Row=47; Col=5; Val='-0.25%'; IsText=$false},
    @{Row=48; Col=4; Val='2.551'; IsText=$false},
    @{Row=48; Col=5; Val='13.84%'; IsText=$false},
    @{Row=49; Col=4; Val='0.001995'; IsText=$false},
    @{Row=49; Col=5; Val='-31.19%'; IsText=$false},
    @{Row=50; Col=4; Val='0.00002094'; IsText=$false},
    @{Row=50; Col=5; Val='-0.25%'; IsText=$false},
    @{Row=51; Col=4; Val='0.0001995'; IsText=$false},
    @{Row=51; Col=5; Val='-0.25%'; IsText=$false}
)

foreach ($u in $updates) {
    if ($u.IsText) {
        $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val
    } else {
        Set-TextValue $ws $u.Row $u.Col $u.Val
    }
}
